# Add data for 2025-02-06
# Update violent crime year-to-date figures across Citywide Totals,
# By Neighborhood, and per-neighborhood sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 560
$ws.Range("L3").Value = 555
$ws.Range("D4").Value = 262
$ws.Range("F4").Value = 235
$ws.Range("I4").Value = 162
$ws.Range("L4").Value = 150
$ws.Range("L5").Value = 46
$ws.Range("K6").Value = 857
$ws.Range("L6").Value = 631
$ws.Range("D7").Value = 2659
$ws.Range("F7").Value = 2079
$ws.Range("I7").Value = 2127
$ws.Range("K7").Value = 2312
$ws.Range("L7").Value = 1942

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L7").Value = 71
$ws.Range("L8").Value = 103
$ws.Range("L11").Value = 32
$ws.Range("L14").Value = 8
$ws.Range("L18").Value = 14
$ws.Range("L19").Value = 68
$ws.Range("L23").Value = 18
$ws.Range("L25").Value = 11
$ws.Range("L27").Value = 19
$ws.Range("L29").Value = 103
$ws.Range("L33").Value = 81
$ws.Range("L36").Value = 32
$ws.Range("K37").Value = 72
$ws.Range("L37").Value = 62
$ws.Range("L41").Value = 10
$ws.Range("L42").Value = 70
$ws.Range("L47").Value = 16
$ws.Range("L52").Value = 37
$ws.Range("L54").Value = 36
$ws.Range("D63").Value = 72
$ws.Range("F63").Value = 38
$ws.Range("I63").Value = 29
$ws.Range("L63").Value = 10
$ws.Range("L65").Value = 38
$ws.Range("L67").Value = 58
$ws.Range("L77").Value = 12
$ws.Range("L79").Value = 53
$ws.Range("L83").Value = 39
$ws.Range("L85").Value = 97
$ws.Range("L86").Value = 15
$ws.Range("L88").Value = 31
$ws.Range("L91").Value = 30
$ws.Range("L99").Value = 35
$ws.Range("D101").Value = 2659
$ws.Range("F101").Value = 2079
$ws.Range("I101").Value = 2127
$ws.Range("K101").Value = 2312
$ws.Range("L101").Value = 1942

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L3").Value = 3
$ws.Range("L7").Value = 8

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 23
$ws.Range("L7").Value = 71

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 24
$ws.Range("L6").Value = 20
$ws.Range("L7").Value = 97

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("L6").Value = 9
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L6").Value = 10
$ws.Range("L7").Value = 39

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L3").Value = 31
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L3").Value = 19
$ws.Range("K6").Value = 24
$ws.Range("L6").Value = 19
$ws.Range("K7").Value = 72
$ws.Range("L7").Value = 62

$ws = $wb.Worksheets.Item("New City")
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 38

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 10
$ws.Range("L7").Value = 35

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L2").Value = 17
$ws.Range("L7").Value = 58

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 7
$ws.Range("L6").Value = 18
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 34
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 25
$ws.Range("L3").Value = 16
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 68

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("L6").Value = 2
$ws.Range("L7").Value = 10

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L6").Value = 33
$ws.Range("L7").Value = 70

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L2").Value = 3
$ws.Range("L7").Value = 18

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 30

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 21
$ws.Range("L3").Value = 16
$ws.Range("L7").Value = 53

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 14

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 15
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L3").Value = 6
$ws.Range("L7").Value = 11

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L2").Value = 9
$ws.Range("L7").Value = 16

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("L2").Value = 7
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 31

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L3").Value = 5
$ws.Range("L4").Value = 3
$ws.Range("L7").Value = 19

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L4").Value = 12
$ws.Range("L6").Value = 15

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L6").Value = 3
$ws.Range("L7").Value = 12
